# Modify max stock for HGTs (Heavy Goods Truck_Maximum stock), sheet "UCT1"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")

# --- Row 28: UC_HGT_Stock year 2018 row. Multiplier 1.1 -> 1.16 ---
$ws.Range("H28").Formula = "=H45*1.16"
$ws.Range("I28:AH28").Formula = "=I45*1.16"

# --- Row 29: year 2022 row. Drop the extra *1.05 growth, just mirror H28 ---
$ws.Range("H29").Formula = "=H28"
# I29:AH29 already reference the row above (=I28 etc.) and recalc automatically.

# --- Row 30: year 2030 row. Multiplier on I:AH columns 1.3 -> 1.25 (H30 already *1.25) ---
$ws.Range("I30:AH30").Formula = "=I29*1.25"

# --- Match number format of H28:H31 to the rest of the row (I:AH), i.e. whole-number "0" ---
$ws.Range("H28:H31").NumberFormat = "0"

# --- Row 82: stray leftover cells no longer used - clear them out ---
$ws.Range("G82:H82").ClearContents()

# --- Restore the view/selection state on the UCT1 sheet ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("K37").Select()
